# Alumni-Map-Project-Mockup.pptx — "Reparsed the data, added documentation"
#
# Appends a new Title-Only slide at the end of the deck with a note
# describing the next round of project tasks (Hv plot, readme, more
# figures, mines-vs-gp-dept population comparison).

$p = $ppt.ActivePresentation

# Add a new slide after the last one, using the "Title Only" layout
# (same CustomLayout/slideLayout6.xml used by all the other content
# slides in this deck). ppLayoutTitleOnly = 11.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 11)

$titleShape = $s.Shapes.Item(1)
$tf = $titleShape.TextFrame
$tr = $tf.TextRange

# Full note text.
$tr.Text = "Hv plot, flush out readme, more project tasks, draw different plots and describe in the planned methodology, 4-5 figs, comparing mines pop to gp dept pop"

# "gp" (in "... comparing mines pop to gp dept pop") was flagged by the
# spell checker, which splits it into its own run in the source deck.
$fullText = $tr.Text
$gpStart = $fullText.IndexOf("gp dept pop") + 1
$gpRange = $tr.Characters($gpStart, 2)
$gpRange.Text = "gp"

# The long title text overflows the default title box, so PowerPoint grew
# the box and shrank the text to fit (normAutofit).
$titleShape.Height = 158.597
$tf.AutoSize = 2
